$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds the single data record in this sample sheet. The dataset is
# being "balanced" by swapping several of the record's values for new
# numbers. Every value in row 2 is stored as text (a shared string) rather
# than a numeric cell, so each target cell is briefly marked as Text before
# the new value is written (otherwise a numeric-looking string like "10000"
# would be auto-converted into a real number). The Text formatting is then
# removed again so the cells end up as plain, unstyled text cells - matching
# how the rest of the row is already stored.
$changes = @{
    "A2" = "10000"
    "E2" = "5"
    "F2" = "0"
    "G2" = "0"
    "H2" = "50"
    "K2" = "9"
    "M2" = "500"
}

foreach ($addr in $changes.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$addr]
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Style = "Normal"
}
